{"js": "// Update the date line and the 25 \"three-digit \u00f7 one-digit\" answer cells\n// to the new values, preserving each run's existing formatting.\nconst replacements = [\n  [\"2024-02-27 Tuesday\", \"2024-02-28 Wednesday\"],\n  [\"431\u00f77=61, 4\", \"804\u00f73=268, 0\"],\n  [\"395\u00f74=98, 3\", \"232\u00f74=58, 0\"],\n  [\"706\u00f74=176, 2\", \"473\u00f79=52, 5\"],\n  [\"335\u00f76=55, 5\", \"482\u00f78=60, 2\"],\n  [\"196\u00f73=65, 1\", \"400\u00f75=80, 0\"],\n  [\"874\u00f74=218, 2\", \"225\u00f78=28, 1\"],\n  [\"159\u00f79=17, 6\", \"974\u00f74=243, 2\"],\n  [\"260\u00f74=65, 0\", \"612\u00f72=306, 0\"],\n  [\"373\u00f72=186, 1\", \"562\u00f73=187, 1\"],\n  [\"530\u00f72=265, 0\", \"900\u00f75=180, 0\"],\n  [\"728\u00f76=121, 2\", \"273\u00f77=39, 0\"],\n  [\"693\u00f73=231, 0\", \"485\u00f79=53, 8\"],\n  [\"151\u00f73=50, 1\", \"956\u00f78=119, 4\"],\n  [\"228\u00f77=32, 4\", \"581\u00f75=116, 1\"],\n  [\"499\u00f72=249, 1\", \"403\u00f78=50, 3\"],\n  [\"799\u00f72=399, 1\", \"995\u00f77=142, 1\"],\n  [\"348\u00f72=174, 0\", \"968\u00f79=107, 5\"],\n  [\"167\u00f73=55, 2\", \"921\u00f77=131, 4\"],\n  [\"174\u00f76=29, 0\", \"457\u00f78=57, 1\"],\n  [\"784\u00f78=98, 0\", \"820\u00f72=410, 0\"],\n  [\"668\u00f73=222, 2\", \"420\u00f76=70, 0\"],\n  [\"449\u00f73=149, 2\", \"214\u00f79=23, 7\"],\n  [\"276\u00f79=30, 6\", \"965\u00f79=107, 2\"],\n  [\"429\u00f72=214, 1\", \"679\u00f72=339, 1\"],\n  [\"397\u00f75=79, 2\", \"499\u00f77=71, 2\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 \"three-digit \u00f7 one-digit\" answer cells\n# to the new values, preserving each run's existing formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-27 Tuesday\", \"2024-02-28 Wednesday\"),\n    @(\"431\u00f77=61, 4\", \"804\u00f73=268, 0\"),\n    @(\"395\u00f74=98, 3\", \"232\u00f74=58, 0\"),\n    @(\"706\u00f74=176, 2\", \"473\u00f79=52, 5\"),\n    @(\"335\u00f76=55, 5\", \"482\u00f78=60, 2\"),\n    @(\"196\u00f73=65, 1\", \"400\u00f75=80, 0\"),\n    @(\"874\u00f74=218, 2\", \"225\u00f78=28, 1\"),\n    @(\"159\u00f79=17, 6\", \"974\u00f74=243, 2\"),\n    @(\"260\u00f74=65, 0\", \"612\u00f72=306, 0\"),\n    @(\"373\u00f72=186, 1\", \"562\u00f73=187, 1\"),\n    @(\"530\u00f72=265, 0\", \"900\u00f75=180, 0\"),\n    @(\"728\u00f76=121, 2\", \"273\u00f77=39, 0\"),\n    @(\"693\u00f73=231, 0\", \"485\u00f79=53, 8\"),\n    @(\"151\u00f73=50, 1\", \"956\u00f78=119, 4\"),\n    @(\"228\u00f77=32, 4\", \"581\u00f75=116, 1\"),\n    @(\"499\u00f72=249, 1\", \"403\u00f78=50, 3\"),\n    @(\"799\u00f72=399, 1\", \"995\u00f77=142, 1\"),\n    @(\"348\u00f72=174, 0\", \"968\u00f79=107, 5\"),\n    @(\"167\u00f73=55, 2\", \"921\u00f77=131, 4\"),\n    @(\"174\u00f76=29, 0\", \"457\u00f78=57, 1\"),\n    @(\"784\u00f78=98, 0\", \"820\u00f72=410, 0\"),\n    @(\"668\u00f73=222, 2\", \"420\u00f76=70, 0\"),\n    @(\"449\u00f73=149, 2\", \"214\u00f79=23, 7\"),\n    @(\"276\u00f79=30, 6\", \"965\u00f79=107, 2\"),\n    @(\"429\u00f72=214, 1\", \"679\u00f72=339, 1\"),\n    @(\"397\u00f75=79, 2\", \"499\u00f77=71, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
